# Generate Report for Handoff
#
# The b6e66154-... file is now ready for handoff (status changes from
# "Handed back: in sync with en-US" to "Ready for handoff" on the Overview
# sheet as well as on each language sheet), and the latest handoff
# timestamps for the 65c7733d-... file are refreshed on each language
# sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the b6e66154-... file ---
$wsOverview.Cells.Item(3, 2).Value = "Ready for handoff"   # B3 (zh-cn column)
$wsOverview.Cells.Item(3, 3).Value = "Ready for handoff"   # C3 (de-de column)

# --- zh-cn sheet ---
# Row 2 (65c7733d-... file): refresh the Latest Handoff Datetime
$wsZhCn.Cells.Item(2, 4).Value = "2016-03-09 10:21:39"     # D2
# Row 3 (b6e66154-... file): status + refreshed handoff datetime
$wsZhCn.Cells.Item(3, 2).Value = "Ready for handoff"       # B3
$wsZhCn.Cells.Item(3, 4).Value = "2016-03-09 10:21:39"     # D3

# --- de-de sheet ---
# Row 2 (65c7733d-... file): refresh the Latest Handoff Datetime
$wsDeDe.Cells.Item(2, 4).Value = "2016-03-09 10:21:50"     # D2
# Row 3 (b6e66154-... file): status + refreshed handoff datetime
$wsDeDe.Cells.Item(3, 2).Value = "Ready for handoff"       # B3
$wsDeDe.Cells.Item(3, 4).Value = "2016-03-09 10:21:50"     # D3
